$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Price (column D), Volume/1h (column E)
# $null for Price means the Price cell is unchanged for that row.
$updates = @(
    @{ Row = 2;  D = "57.824.38";  E = "  +1.63%  " },
    @{ Row = 3;  D = "3.053.67";   E = "  +0.46%  " },
    @{ Row = 4;  D = $null;        E = "  +0.04%  " },
    @{ Row = 5;  D = "515.73";     E = "  +0.70%  " },
    @{ Row = 6;  D = "141.41";     E = "  +0.33%  " },
    @{ Row = 7;  D = $null;        E = "  +0.00%  " },
    @{ Row = 8;  D = "0.436";      E = "  +0.80%  " },
    @{ Row = 9;  D = "7.28";       E = "  +1.66%  " },
    @{ Row = 10; D = $null;        E = "  -1.18%  " },
    @{ Row = 11; D = $null;        E = "  +1.08%  " },
    @{ Row = 12; D = "3.575.57";   E = "  +0.74%  " },
    @{ Row = 13; D = $null;        E = "  +2.91%  " },
    @{ Row = 14; D = "26.15";      E = "  +2.84%  " },
    @{ Row = 15; D = "0.0000163";  E = "  -0.39%  " },
    @{ Row = 16; D = "57.852.30";  E = "  +1.81%  " },
    @{ Row = 17; D = "3.048.93";   E = "  +0.56%  " },
    @{ Row = 18; D = $null;        E = "  +2.32%  " },
    @{ Row = 19; D = "12.79";      E = "  -3.05%  " },
    @{ Row = 20; D = "8.05";       E = "  -0.53%  " },
    @{ Row = 21; D = "330.70";     E = "  -1.29%  " },
    @{ Row = 22; D = "1.00";       E = "  +0.00%  " },
    @{ Row = 23; D = "0.500";      E = "  -0.58%  " },
    @{ Row = 24; D = "65.36";      E = "  +0.87%  " },
    @{ Row = 25; D = "0.170";      E = "  +2.10%  " },
    @{ Row = 26; D = $null;        E = "  -0.02%  " },
    @{ Row = 27; D = ("0.0{0}0896" -f [char]0x2083); E = "  -4.86%  " },
    @{ Row = 28; D = "6.41";       E = "  -0.72%  " },
    @{ Row = 29; D = "7.23";       E = "  +6.47%  " },
    @{ Row = 30; D = $null;        E = "  +1.10%  " },
    @{ Row = 31; D = $null;        E = "  +1.66%  " },
    @{ Row = 32; D = "20.64";      E = "  +0.63%  " },
    @{ Row = 33; D = "154.60";     E = "  +1.08%  " },
    @{ Row = 34; D = "4.51";       E = "  +0.07%  " },
    @{ Row = 35; D = $null;        E = "  +2.58%  " },
    @{ Row = 36; D = "26.91";      E = "  -0.95%  " },
    @{ Row = 37; D = "1.26";       E = "  +2.50%  " },
    @{ Row = 38; D = "0.0680";     E = "  +2.16%  " },
    @{ Row = 39; D = "3.092.28";   E = "  +0.59%  " },
    @{ Row = 40; D = "3.91";       E = "  +2.46%  " },
    @{ Row = 41; D = "36.56";      E = "  -0.22%  " },
    @{ Row = 42; D = $null;        E = "  +0.10%  " },
    @{ Row = 43; D = "0.652";      E = "  -1.37%  " },
    @{ Row = 44; D = "2.270.19";   E = "  +2.64%  " },
    @{ Row = 45; D = $null;        E = "  +5.26%  " },
    @{ Row = 46; D = $null;        E = "  +0.60%  " },
    @{ Row = 47; D = "20.63";      E = "  +3.83%  " },
    @{ Row = 48; D = "5.92";       E = "  +0.81%  " },
    @{ Row = 49; D = "0.935";      E = "  -0.65%  " },
    @{ Row = 50; D = "0.740";      E = "  +9.51%  " },
    @{ Row = 51; D = "252.53";     E = "  +8.91%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel stores the value as literal text
        # rather than re-interpreting it as a number (this mirrors how a user
        # would force-enter text like "1.00" or "330.70" into a cell).
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Value = "'" + $u.D
        # Restore the default "Normal" style so no extra number-format /
        # quote-prefix styling is left behind on the cell.
        $cell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
